# Insert a new weekly record at row 127 for "Vega Modelo de Temuco" / Coliflor.
# All subsequent rows (old 127..246) shift down by one (new 128..247),
# and the sheet's used range grows from A1:R246 to A1:R247.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 127 (and everything below it) down one row.
$ws.Rows("127").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A127").Value = 10
$ws.Range("B127").Value = "Vega Modelo de Temuco"
$ws.Range("C127").Value = "La Araucanía"
$ws.Range("D127").Value = 44484
$ws.Range("E127").Value = 9
$ws.Range("F127").Value = 100112008
$ws.Range("G127").Value = "Coliflor"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 600
$ws.Range("K127").Value = 700
$ws.Range("L127").Value = 800
$ws.Range("M127").Value = 750
$ws.Range("N127").Value = "`$/unidad"
$ws.Range("O127").Value = "Región Metropolitana"
$ws.Range("P127").Value = 750
$ws.Range("Q127").Value = 1
$ws.Range("R127").Value = "Hortaliza"
